$p = $ppt.ActivePresentation

# Add a new "Title and Content" slide (layout 2) at the end of the deck.
$idx = $p.Slides.Count + 1
$s = $p.Slides.Add($idx, 2)

# Content placeholder: 5 paragraphs, 3rd is blank, last one indented one level.
$content = $s.Shapes.Item(2).TextFrame.TextRange
$content.Text = "1 container for back end"
$content.InsertAfter("`rMultiple containers for front end.") | Out-Null
$content.InsertAfter("`r") | Out-Null
$content.InsertAfter("`rNo persistency for back end if containerized") | Out-Null
$content.InsertAfter("`rWould require redesigning.") | Out-Null

# "Would require redesigning." is indented to the second outline level.
$content.Paragraphs(5).IndentLevel = 2
